$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Non insulated v4 TUV")
$ws.Name = "Data"
